# Horarios actualizados Linea 141 - 370
# Applies the 06:43:12 scrape update across all three sheets:
#   - LP1912: rows 41-60 updated/appended (48 -> 55 data rows)
#   - LP1912-215: rows 15-16 appended (9 -> 11 data rows)
#   - 6203-6173: only the "Ultima actualizacion" timestamp changes

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 06:43:12"
$ws1.Range("A3").Value = "Total filas: 55"

$ws1.Cells.Item(41,1).Value = "06:43:12"
$ws1.Cells.Item(41,2).Value = "07:17"
$ws1.Cells.Item(41,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(41,4).Value = 34
$ws1.Cells.Item(41,5).Value = "LP1912"

$ws1.Cells.Item(42,1).Value = "05:49:10"
$ws1.Cells.Item(42,2).Value = "07:21"
$ws1.Cells.Item(42,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(42,4).Value = 92
$ws1.Cells.Item(42,5).Value = "LP1912"

$ws1.Cells.Item(43,1).Value = "06:14:19"
$ws1.Cells.Item(43,2).Value = "07:23"
$ws1.Cells.Item(43,3).Value = "10_OLMOS"
$ws1.Cells.Item(43,4).Value = 69
$ws1.Cells.Item(43,5).Value = "LP1912"

$ws1.Cells.Item(44,1).Value = "06:14:19"
$ws1.Cells.Item(44,2).Value = "07:31"
$ws1.Cells.Item(44,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(44,4).Value = 77
$ws1.Cells.Item(44,5).Value = "LP1912"

$ws1.Cells.Item(45,1).Value = "05:49:10"
$ws1.Cells.Item(45,2).Value = "07:32"
$ws1.Cells.Item(45,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(45,4).Value = 103
$ws1.Cells.Item(45,5).Value = "LP1912"

$ws1.Cells.Item(46,1).Value = "05:49:10"
$ws1.Cells.Item(46,2).Value = "07:32"
$ws1.Cells.Item(46,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(46,4).Value = 103
$ws1.Cells.Item(46,5).Value = "LP1912"

$ws1.Cells.Item(47,1).Value = "05:49:10"
$ws1.Cells.Item(47,2).Value = "07:32"
$ws1.Cells.Item(47,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(47,4).Value = 103
$ws1.Cells.Item(47,5).Value = "LP1912"

$ws1.Cells.Item(48,1).Value = "05:49:10"
$ws1.Cells.Item(48,2).Value = "07:37"
$ws1.Cells.Item(48,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(48,4).Value = 108
$ws1.Cells.Item(48,5).Value = "LP1912"

$ws1.Cells.Item(49,1).Value = "05:49:10"
$ws1.Cells.Item(49,2).Value = "07:39"
$ws1.Cells.Item(49,3).Value = "10_OLMOS"
$ws1.Cells.Item(49,4).Value = 110
$ws1.Cells.Item(49,5).Value = "LP1912"

$ws1.Cells.Item(50,1).Value = "06:14:19"
$ws1.Cells.Item(50,2).Value = "07:47"
$ws1.Cells.Item(50,3).Value = "14_ABASTO"
$ws1.Cells.Item(50,4).Value = 93
$ws1.Cells.Item(50,5).Value = "LP1912"

$ws1.Cells.Item(51,1).Value = "05:49:10"
$ws1.Cells.Item(51,2).Value = "07:48"
$ws1.Cells.Item(51,3).Value = "14_ABASTO"
$ws1.Cells.Item(51,4).Value = 119
$ws1.Cells.Item(51,5).Value = "LP1912"

$ws1.Cells.Item(52,1).Value = "06:14:19"
$ws1.Cells.Item(52,2).Value = "07:51"
$ws1.Cells.Item(52,3).Value = "215D_EL PATO"
$ws1.Cells.Item(52,4).Value = 97
$ws1.Cells.Item(52,5).Value = "LP1912"

$ws1.Cells.Item(53,1).Value = "06:43:12"
$ws1.Cells.Item(53,2).Value = "07:52"
$ws1.Cells.Item(53,3).Value = "215D_EL PATO"
$ws1.Cells.Item(53,4).Value = 69
$ws1.Cells.Item(53,5).Value = "LP1912"

$ws1.Cells.Item(54,1).Value = "06:14:19"
$ws1.Cells.Item(54,2).Value = "08:00"
$ws1.Cells.Item(54,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(54,4).Value = 106
$ws1.Cells.Item(54,5).Value = "LP1912"

$ws1.Cells.Item(55,1).Value = "06:43:12"
$ws1.Cells.Item(55,2).Value = "08:03"
$ws1.Cells.Item(55,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(55,4).Value = 80
$ws1.Cells.Item(55,5).Value = "LP1912"

$ws1.Cells.Item(56,1).Value = "06:14:19"
$ws1.Cells.Item(56,2).Value = "08:12"
$ws1.Cells.Item(56,3).Value = "15_ABASTO"
$ws1.Cells.Item(56,4).Value = 118
$ws1.Cells.Item(56,5).Value = "LP1912"

$ws1.Cells.Item(57,1).Value = "06:43:12"
$ws1.Cells.Item(57,2).Value = "08:21"
$ws1.Cells.Item(57,3).Value = "26_HERNANDEZ"
$ws1.Cells.Item(57,4).Value = 98
$ws1.Cells.Item(57,5).Value = "LP1912"

$ws1.Cells.Item(58,1).Value = "06:43:12"
$ws1.Cells.Item(58,2).Value = "08:23"
$ws1.Cells.Item(58,3).Value = "215B_EL PATO"
$ws1.Cells.Item(58,4).Value = 100
$ws1.Cells.Item(58,5).Value = "LP1912"

$ws1.Cells.Item(59,1).Value = "06:43:12"
$ws1.Cells.Item(59,2).Value = "08:23"
$ws1.Cells.Item(59,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(59,4).Value = 100
$ws1.Cells.Item(59,5).Value = "LP1912"

$ws1.Cells.Item(60,1).Value = "06:43:12"
$ws1.Cells.Item(60,2).Value = "08:27"
$ws1.Cells.Item(60,3).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(60,4).Value = 104
$ws1.Cells.Item(60,5).Value = "LP1912"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 06:43:12"
$ws2.Range("A3").Value = "Total filas: 11"

$ws2.Cells.Item(15,1).Value = "06:43:12"
$ws2.Cells.Item(15,2).Value = "07:52"
$ws2.Cells.Item(15,3).Value = "215D_EL PATO"
$ws2.Cells.Item(15,4).Value = 69
$ws2.Cells.Item(15,5).Value = "LP1912"

$ws2.Cells.Item(16,1).Value = "06:43:12"
$ws2.Cells.Item(16,2).Value = "08:23"
$ws2.Cells.Item(16,3).Value = "215B_EL PATO"
$ws2.Cells.Item(16,4).Value = 100
$ws2.Cells.Item(16,5).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 06:43:12"

